$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "survey" sheet: add a "sex" question (with its own explanatory note)
#    right before the weight-for-age plot's custom template note.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Insert a new row above the existing "age" question (row 46) for a note
# explaining the upcoming weight-for-age plot questions.
$survey.Rows.Item(46).Insert()
$survey.Range("A46").Value = "note"
$survey.Range("D46").Value = "The following data will be used to generate a weight for age plot."

# The "age" question is now on row 47 - clarify the label and add a
# constraint message.
$survey.Range("D47").Value = "Enter age (in years):"
$survey.Range("E47").Value = "Must be less than 20."

# The "weight" question is now on row 48 - clarify the label.
$survey.Range("D48").Value = "Enter weight (in lbs):"

# Insert a new row after "weight" (row 49) for a "sex" select-one question.
$survey.Rows.Item(49).Insert()
$survey.Range("A49").Value = "select_one sexes"
$survey.Range("C49").Value = "sex"
$survey.Range("D49").Value = "Enter sex:"

# ---------------------------------------------------------------------------
# 2. "choices" sheet: add the "sexes" choice list (male / female) used by
#    the new "select_one sexes" question above.
# ---------------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("A23").Value = "sexes"
$choices.Range("B23").Value = "male"
$choices.Range("C23").Value = "male"

$choices.Range("A24").Value = "sexes"
$choices.Range("B24").Value = "female"
$choices.Range("C24").Value = "female"
